$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add new "Revise" column (G) ------------------------------------
# Reuse the *old* F-column formats for G (G keeps the full-box border that
# F used to have) by copying F1:F11 formats across before F's own border
# is altered.
$ws.Range("F1:F11").Copy()
$ws.Range("G1:G11").PasteSpecial(-4122)   # xlPasteFormats

# Header
$ws.Range("G1").Value = "Revise"

# Per-row verdicts for the new "Revise" column
$ws.Range("G2").Value  = "Pass"
$ws.Range("G3").Value  = "Pass"
$ws.Range("G4").Value  = "Failed"
$ws.Range("G5").Value  = "Pass"
$ws.Range("G6").Value  = "Pass"
$ws.Range("G7").Value  = "Pass"
$ws.Range("G8").Value  = "Pass"
$ws.Range("G9").Value  = "Failed"
$ws.Range("G10").Value = "Pass"
$ws.Range("G11").Value = "Failed"

# --- 2. F column no longer is the right-most column: drop its right edge
$fBorder = $ws.Range("F1:F11").Borders.Item(10)   # xlEdgeRight
$fBorder.LineStyle = -4142                         # xlLineStyleNone

# --- 3. Update the active selection/cursor position ---------------------
$ws.Range("H10").Select()
